# Update GPapa BOM "In Stock In Bay" (column G) quantities/notes, and
# refresh the selection to G19. Also clears the leftover "Text" style
# (xfId 1 / applyNumberFormat) that was stamped across the data rows,
# matching the author's formatting cleanup.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the stray per-cell style (s="1") from the whole data range so the
# cells fall back to the default/general style, same as the committed file.
$ws.Range("A2:G56").ClearFormats()

# Update in-stock notes for the capacitor / LED rows.
$ws.Range("G8").Value = "43 in the bay"
$ws.Range("G9").Value = "100+ in the bay"
$ws.Range("G10").Value = 0
$ws.Range("G11").Value = "50~ in the Bay"
$ws.Range("G12").Value = "25 in the Bay"
$ws.Range("G13").Value = 0
$ws.Range("G14").Value = "100+ in the Bay"
$ws.Range("G15").Value = "49 in the bay"
$ws.Range("G16").Value = "100+ in bay"
$ws.Range("G17").Value = "100+ in bay"
$ws.Range("G18").Value = "100+ in bay"

# Move the cursor/selection to G19, matching the saved view state.
$ws.Range("G19").Select() | Out-Null
